# Add "label_for_cases_en" / "label_for_cases_fra" columns to the
# Modules_and_forms sheet, inserted right after "default_fra" and before
# "icon_filepath". Existing columns to the right shift over by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Modules_and_forms")

# Insert two new blank columns at E:F, shifting icon_filepath,
# audio_filepath, unique_id from E:G to G:I.
$ws.Range("E1:F3").EntireColumn.Insert()

# New header cells
$ws.Range("E1").Value = "label_for_cases_en"
$ws.Range("F1").Value = "label_for_cases_fra"

# Row 2 (Module) gets "Cases" for both new columns
$ws.Range("E2").Value = "Cases"
$ws.Range("F2").Value = "Cases"

# Row 3 (Form) leaves the new columns blank (no assignment needed -
# the inserted cells start out empty).
